# Update the lattice multiplication exercises table: swap each
# cell's exercise ("A x B", the spaced-out digits of B, the divider,
# and the two partial-product placeholder rows) for a new one.
#
# NOTE: this runtime's Range.Find always scans from the start of the
# document regardless of which Range object it is invoked on, so a
# Find scoped to a single cell can still (incorrectly) edit a match
# in an earlier cell when the search text repeats across cells (e.g.
# "2|    |" shows up in many cells). Table cell content here is short
# and entirely replaced anyway, so we sidestep Find altogether and set
# each target cell's Range.Text directly via Table.Cell(row, col) --
# that indexing is reliable and keeps every write inside its own cell.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$nl = [char]11   # vertical-tab == a w:br line break inside the run

$t.Cell(1,1).Range.Text = "62 x 13" + $nl + "  1    3" + $nl + "  ----" + $nl + "6|    |" + $nl + "2|    |"
$t.Cell(1,2).Range.Text = "54 x 55" + $nl + "  5    5" + $nl + "  ----" + $nl + "5|    |" + $nl + "4|    |"
$t.Cell(1,3).Range.Text = "54 x 51" + $nl + "  5    1" + $nl + "  ----" + $nl + "5|    |" + $nl + "4|    |"
$t.Cell(2,1).Range.Text = "15 x 55" + $nl + "  5    5" + $nl + "  ----" + $nl + "1|    |" + $nl + "5|    |"
$t.Cell(2,2).Range.Text = "71 x 19" + $nl + "  1    9" + $nl + "  ----" + $nl + "7|    |" + $nl + "1|    |"
$t.Cell(2,3).Range.Text = "89 x 90" + $nl + "  9    0" + $nl + "  ----" + $nl + "8|    |" + $nl + "9|    |"
$t.Cell(3,1).Range.Text = "80 x 75" + $nl + "  7    5" + $nl + "  ----" + $nl + "8|    |" + $nl + "0|    |"
$t.Cell(3,2).Range.Text = "12 x 95" + $nl + "  9    5" + $nl + "  ----" + $nl + "1|    |" + $nl + "2|    |"
$t.Cell(3,3).Range.Text = "69 x 57" + $nl + "  5    7" + $nl + "  ----" + $nl + "6|    |" + $nl + "9|    |"
$t.Cell(4,1).Range.Text = "17 x 19" + $nl + "  1    9" + $nl + "  ----" + $nl + "1|    |" + $nl + "7|    |"
$t.Cell(4,2).Range.Text = "19 x 78" + $nl + "  7    8" + $nl + "  ----" + $nl + "1|    |" + $nl + "9|    |"
$t.Cell(4,3).Range.Text = "19 x 33" + $nl + "  3    3" + $nl + "  ----" + $nl + "1|    |" + $nl + "9|    |"
$t.Cell(5,1).Range.Text = "89 x 98" + $nl + "  9    8" + $nl + "  ----" + $nl + "8|    |" + $nl + "9|    |"
$t.Cell(5,2).Range.Text = "41 x 74" + $nl + "  7    4" + $nl + "  ----" + $nl + "4|    |" + $nl + "1|    |"
$t.Cell(5,3).Range.Text = "10 x 23" + $nl + "  2    3" + $nl + "  ----" + $nl + "1|    |" + $nl + "0|    |"
